$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '35.454.83'
$ws.Range('E2').Value = '  +2.10%  '

# Row 3
$ws.Range('D3').Value = '1.899.92'
$ws.Range('E3').Value = '  +2.23%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '247.96'
$ws.Range('E5').Value = '  +1.64%  '

# Row 6
$ws.Range('E6').Value = '  +3.25%  '

# Row 7
$ws.Range('E7').Value = '  +0.00%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '43.19'
$ws.Range('E8').Value = '  +3.67%  '

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.359'
$ws.Range('E9').Value = '  +6.06%  '

# Row 10
$ws.Range('E10').Value = '  +9.16%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0748'
$ws.Range('E11').Value = '  +2.83%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0979'
$ws.Range('E12').Value = '  +1.48%  '

# Row 13
$ws.Range('E13').Value = '  +9.73%  '

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.789'
$ws.Range('E14').Value = '  +12.16%  '

# Row 15
$ws.Range('D15').Value = '2.171.40'
$ws.Range('E15').Value = '  +1.92%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '5.00'
$ws.Range('E16').Value = '  +4.63%  '

# Row 17
$ws.Range('D17').Value = '1.909.25'
$ws.Range('E17').Value = '  +2.52%  '

# Row 18
$ws.Range('D18').Value = '35.465.70'
$ws.Range('E18').Value = '  +2.22%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '73.67'
$ws.Range('E19').Value = '  +2.54%  '

# Row 20
$ws.Range('D20').Value = '0.0₃0829'
$ws.Range('E20').Value = '  +2.92%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '245.05'
$ws.Range('E21').Value = '  +1.34%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '12.91'
$ws.Range('E22').Value = '  +3.63%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.25'
$ws.Range('E23').Value = '  +8.76%  '

# Row 24
$ws.Range('E24').Value = '  +9.04%  '

# Row 25
$ws.Range('E25').Value = '  -0.06%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.18'
$ws.Range('E26').Value = '  +2.06%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '166.98'
$ws.Range('E27').Value = '  +2.58%  '

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '8.57'
$ws.Range('E28').Value = '  +3.61%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '18.38'
$ws.Range('E29').Value = '  +2.41%  '

# Row 30
$ws.Range('E30').Value = '  +2.17%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '4.36'
$ws.Range('E31').Value = '  +5.65%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.0601'
$ws.Range('E32').Value = '  +6.02%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.23'
$ws.Range('E33').Value = '  +3.70%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.86'
$ws.Range('E34').Value = '  +18.18%  '

# Row 35
$ws.Range('E35').Value = '  -0.09%  '

# Row 36
$ws.Range('E36').Value = '  -13.29%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.858'
$ws.Range('E37').Value = '  +4.74%  '

# Row 38
$ws.Range('E38').Value = '  +2.65%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0711'
$ws.Range('E39').Value = '  +7.57%  '

# Row 40
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0222'
$ws.Range('E40').Value = '  +6.65%  '

# Row 41
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '99.80'
$ws.Range('E41').Value = '  +2.95%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '17.14'
$ws.Range('E42').Value = '  +2.27%  '

# Row 43
$ws.Range('E43').Value = '  +2.79%  '

# Row 44
$ws.Range('D44').Value = '1.337.66'
$ws.Range('E44').Value = '  +4.84%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '13.54'
$ws.Range('E45').Value = '  +15.42%  '

# Row 46
$ws.Range('E46').Value = '  +4.50%  '

# Row 47
$ws.Range('E47').Value = '  -2.85%  '

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.42'
$ws.Range('E48').Value = '  +0.92%  '

# Row 49
$ws.Range('E49').Value = '  +0.67%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '6.38'
$ws.Range('E50').Value = '  +2.65%  '

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '42.51'
$ws.Range('E51').Value = '  +1.07%  '
